$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 (subject id 10)
$row12 = @("G", 4, 2, 3, 2, 2, 2, 3, 4, 3, 2, 3, 4, 3, 1, 5)
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, 2 + $i).Value = $row12[$i]
}

# Row 13 (subject id 11)
$row13 = @("F", 5, 5, 4, 3, 5, 4, 3, 6, 5, 2, 5, 5, 6, 2, 5)
for ($i = 0; $i -lt $row13.Length; $i++) {
    $ws.Cells.Item(13, 2 + $i).Value = $row13[$i]
}

# Update the selection to match the recorded active cell after the edits
$ws.Range("Q13").Select()
